$wb = $excel.ActiveWorkbook

# Use an existing sheet as a style/format template, then clear its content so the new
# sheet picks up the same sheetFormatPr / base formatting as the other sheets in the book.
$headerSrc = $wb.Worksheets.Item("Solver Team Data")
$tmpl = $wb.Worksheets.Item("Partner Solver Weights")
$tmpl.Copy($null, $tmpl)
$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "Partner Match"
$ws.Cells.Clear()

# Header row
$ws.Cells.Item(1,1).Value = "Partners"
$ws.Cells.Item(1,2).Value = "Solvers"
$ws.Cells.Item(1,3).Value = "Count"
$ws.Cells.Item(1,4).Value = "Comments"

# Copy the header formatting (bold font, thin border, centered) from an existing header row
$headerSrc.Range("A1:D1").Copy()
$ws.Range("A1:D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data rows: one row per partner, tracking linked solvers/count/comments
$ws.Cells.Item(2,1).Value = "Access Afya"
$ws.Cells.Item(2,2).Value = "None"
$ws.Cells.Item(2,3).Value = 0
$ws.Cells.Item(2,4).Value = "None"
$ws.Cells.Item(3,1).Value = "American Family Insurance Institute for Corporate and Social Impact"
$ws.Cells.Item(3,2).Value = "None"
$ws.Cells.Item(3,3).Value = 0
$ws.Cells.Item(3,4).Value = "None"
$ws.Cells.Item(4,1).Value = "American Student Assistance (ASA)"
$ws.Cells.Item(4,2).Value = "None"
$ws.Cells.Item(4,3).Value = 0
$ws.Cells.Item(4,4).Value = "None"
$ws.Cells.Item(5,1).Value = "Americares"
$ws.Cells.Item(5,2).Value = "None"
$ws.Cells.Item(5,3).Value = 0
$ws.Cells.Item(5,4).Value = "None"
$ws.Cells.Item(6,1).Value = "Antropia ESSEC"
$ws.Cells.Item(6,2).Value = "None"
$ws.Cells.Item(6,3).Value = 0
$ws.Cells.Item(6,4).Value = "None"
$ws.Cells.Item(7,1).Value = "AutoCognita"
$ws.Cells.Item(7,2).Value = "None"
$ws.Cells.Item(7,3).Value = 0
$ws.Cells.Item(7,4).Value = "None"
$ws.Cells.Item(8,1).Value = "Best Buy"
$ws.Cells.Item(8,2).Value = "None"
$ws.Cells.Item(8,3).Value = 0
$ws.Cells.Item(8,4).Value = "None"
$ws.Cells.Item(9,1).Value = "Blue Haven Initiative"
$ws.Cells.Item(9,2).Value = "None"
$ws.Cells.Item(9,3).Value = 0
$ws.Cells.Item(9,4).Value = "None"
$ws.Cells.Item(10,1).Value = "BMW Foundation Herbert Quandt"
$ws.Cells.Item(10,2).Value = "None"
$ws.Cells.Item(10,3).Value = 0
$ws.Cells.Item(10,4).Value = "None"
$ws.Cells.Item(11,1).Value = "Cambridge Associates"
$ws.Cells.Item(11,2).Value = "None"
$ws.Cells.Item(11,3).Value = 0
$ws.Cells.Item(11,4).Value = "None"
$ws.Cells.Item(12,1).Value = "Capital One"
$ws.Cells.Item(12,2).Value = "None"
$ws.Cells.Item(12,3).Value = 0
$ws.Cells.Item(12,4).Value = "None"
$ws.Cells.Item(13,1).Value = "Care 2 Communities"
$ws.Cells.Item(13,2).Value = "None"
$ws.Cells.Item(13,3).Value = 0
$ws.Cells.Item(13,4).Value = "None"
$ws.Cells.Item(14,1).Value = "Cast Collective"
$ws.Cells.Item(14,2).Value = "None"
$ws.Cells.Item(14,3).Value = 0
$ws.Cells.Item(14,4).Value = "None"
$ws.Cells.Item(15,1).Value = "Clint Taylor"
$ws.Cells.Item(15,2).Value = "None"
$ws.Cells.Item(15,3).Value = 0
$ws.Cells.Item(15,4).Value = "None"
$ws.Cells.Item(16,1).Value = "Clorox"
$ws.Cells.Item(16,2).Value = "None"
$ws.Cells.Item(16,3).Value = 0
$ws.Cells.Item(16,4).Value = "None"
$ws.Cells.Item(17,1).Value = "Closed Loop Partners"
$ws.Cells.Item(17,2).Value = "None"
$ws.Cells.Item(17,3).Value = 0
$ws.Cells.Item(17,4).Value = "None"
$ws.Cells.Item(18,1).Value = "Comcast NBCUniversal"
$ws.Cells.Item(18,2).Value = "None"
$ws.Cells.Item(18,3).Value = 0
$ws.Cells.Item(18,4).Value = "None"
$ws.Cells.Item(19,1).Value = "Compassion International"
$ws.Cells.Item(19,2).Value = "None"
$ws.Cells.Item(19,3).Value = 0
$ws.Cells.Item(19,4).Value = "None"
$ws.Cells.Item(20,1).Value = "Conduent"
$ws.Cells.Item(20,2).Value = "None"
$ws.Cells.Item(20,3).Value = 0
$ws.Cells.Item(20,4).Value = "None"
$ws.Cells.Item(21,1).Value = "Covestro"
$ws.Cells.Item(21,2).Value = "None"
$ws.Cells.Item(21,3).Value = 0
$ws.Cells.Item(21,4).Value = "None"
$ws.Cells.Item(22,1).Value = "Danaher"
$ws.Cells.Item(22,2).Value = "None"
$ws.Cells.Item(22,3).Value = 0
$ws.Cells.Item(22,4).Value = "None"
$ws.Cells.Item(23,1).Value = "Deshpande Foundation"
$ws.Cells.Item(23,2).Value = "None"
$ws.Cells.Item(23,3).Value = 0
$ws.Cells.Item(23,4).Value = "None"
$ws.Cells.Item(24,1).Value = "Dubai Cares"
$ws.Cells.Item(24,2).Value = "None"
$ws.Cells.Item(24,3).Value = 0
$ws.Cells.Item(24,4).Value = "None"
$ws.Cells.Item(25,1).Value = "eBay"
$ws.Cells.Item(25,2).Value = "None"
$ws.Cells.Item(25,3).Value = 0
$ws.Cells.Item(25,4).Value = "None"
$ws.Cells.Item(26,1).Value = "EcoAdvisors"
$ws.Cells.Item(26,2).Value = "None"
$ws.Cells.Item(26,3).Value = 0
$ws.Cells.Item(26,4).Value = "None"
$ws.Cells.Item(27,1).Value = "EILEEN FISHER"
$ws.Cells.Item(27,2).Value = "None"
$ws.Cells.Item(27,3).Value = 0
$ws.Cells.Item(27,4).Value = "None"
$ws.Cells.Item(28,1).Value = "Enel Foundation and Strategic Studies Center"
$ws.Cells.Item(28,2).Value = "None"
$ws.Cells.Item(28,3).Value = 0
$ws.Cells.Item(28,4).Value = "None"
$ws.Cells.Item(29,1).Value = "Firefly Innovations"
$ws.Cells.Item(29,2).Value = "None"
$ws.Cells.Item(29,3).Value = 0
$ws.Cells.Item(29,4).Value = "None"
$ws.Cells.Item(30,1).Value = "Georgia-Pacific Foundation"
$ws.Cells.Item(30,2).Value = "None"
$ws.Cells.Item(30,3).Value = 0
$ws.Cells.Item(30,4).Value = "None"
$ws.Cells.Item(31,1).Value = "Gina's Collective"
$ws.Cells.Item(31,2).Value = "None"
$ws.Cells.Item(31,3).Value = 0
$ws.Cells.Item(31,4).Value = "None"
$ws.Cells.Item(32,1).Value = "Global Fund to fight Aids, Tuberculosis and Malaria"
$ws.Cells.Item(32,2).Value = "None"
$ws.Cells.Item(32,3).Value = 0
$ws.Cells.Item(32,4).Value = "None"
$ws.Cells.Item(33,1).Value = "Grupo Salinas"
$ws.Cells.Item(33,2).Value = "None"
$ws.Cells.Item(33,3).Value = 0
$ws.Cells.Item(33,4).Value = "None"
$ws.Cells.Item(34,1).Value = "Henkel"
$ws.Cells.Item(34,2).Value = "None"
$ws.Cells.Item(34,3).Value = 0
$ws.Cells.Item(34,4).Value = "None"
$ws.Cells.Item(35,1).Value = "Ingredion"
$ws.Cells.Item(35,2).Value = "None"
$ws.Cells.Item(35,3).Value = 0
$ws.Cells.Item(35,4).Value = "None"
$ws.Cells.Item(36,1).Value = "Innospark Ventures"
$ws.Cells.Item(36,2).Value = "None"
$ws.Cells.Item(36,3).Value = 0
$ws.Cells.Item(36,4).Value = "None"
$ws.Cells.Item(37,1).Value = "Innovation Norway"
$ws.Cells.Item(37,2).Value = "None"
$ws.Cells.Item(37,3).Value = 0
$ws.Cells.Item(37,4).Value = "None"
$ws.Cells.Item(38,1).Value = "Kevin Przybocki"
$ws.Cells.Item(38,2).Value = "None,AHSA Platform"
$ws.Cells.Item(38,3).Value = 1
$ws.Cells.Item(38,4).Value = "None"
$ws.Cells.Item(39,1).Value = "Klaxoon"
$ws.Cells.Item(39,2).Value = "None"
$ws.Cells.Item(39,3).Value = 0
$ws.Cells.Item(39,4).Value = "None"
$ws.Cells.Item(40,1).Value = "KSF Impact"
$ws.Cells.Item(40,2).Value = "None"
$ws.Cells.Item(40,3).Value = 0
$ws.Cells.Item(40,4).Value = "None"
$ws.Cells.Item(41,1).Value = "Leap Ventures"
$ws.Cells.Item(41,2).Value = "None"
$ws.Cells.Item(41,3).Value = 0
$ws.Cells.Item(41,4).Value = "None"
$ws.Cells.Item(42,1).Value = "Lex Mundi Pro Bono Foundation"
$ws.Cells.Item(42,2).Value = "None"
$ws.Cells.Item(42,3).Value = 0
$ws.Cells.Item(42,4).Value = "None"
$ws.Cells.Item(43,1).Value = "Llamasoft"
$ws.Cells.Item(43,2).Value = "None"
$ws.Cells.Item(43,3).Value = 0
$ws.Cells.Item(43,4).Value = "None"
$ws.Cells.Item(44,1).Value = "Mannin Research"
$ws.Cells.Item(44,2).Value = "None"
$ws.Cells.Item(44,3).Value = 0
$ws.Cells.Item(44,4).Value = "None"
$ws.Cells.Item(45,1).Value = "Merck for Mothers"
$ws.Cells.Item(45,2).Value = "None"
$ws.Cells.Item(45,3).Value = 0
$ws.Cells.Item(45,4).Value = "None"
$ws.Cells.Item(46,1).Value = "Merian Ventures"
$ws.Cells.Item(46,2).Value = "None"
$ws.Cells.Item(46,3).Value = 0
$ws.Cells.Item(46,4).Value = "None"
$ws.Cells.Item(47,1).Value = "MIT Club of Northern California"
$ws.Cells.Item(47,2).Value = "None"
$ws.Cells.Item(47,3).Value = 0
$ws.Cells.Item(47,4).Value = "None"
$ws.Cells.Item(48,1).Value = "Mondi Group"
$ws.Cells.Item(48,2).Value = "None"
$ws.Cells.Item(48,3).Value = 0
$ws.Cells.Item(48,4).Value = "None"
$ws.Cells.Item(49,1).Value = "Morgridge Family Foundation"
$ws.Cells.Item(49,2).Value = "None"
$ws.Cells.Item(49,3).Value = 0
$ws.Cells.Item(49,4).Value = "None"
$ws.Cells.Item(50,1).Value = "National Rongxiang Xu Foundation"
$ws.Cells.Item(50,2).Value = "None"
$ws.Cells.Item(50,3).Value = 0
$ws.Cells.Item(50,4).Value = "None"
$ws.Cells.Item(51,1).Value = "Northrop Grumman Corporation"
$ws.Cells.Item(51,2).Value = "None"
$ws.Cells.Item(51,3).Value = 0
$ws.Cells.Item(51,4).Value = "None"
$ws.Cells.Item(52,1).Value = "Nuvo"
$ws.Cells.Item(52,2).Value = "None"
$ws.Cells.Item(52,3).Value = 0
$ws.Cells.Item(52,4).Value = "None"
$ws.Cells.Item(53,1).Value = "Olam International"
$ws.Cells.Item(53,2).Value = "None"
$ws.Cells.Item(53,3).Value = 0
$ws.Cells.Item(53,4).Value = "None"
$ws.Cells.Item(54,1).Value = "Oliver Wyman Group"
$ws.Cells.Item(54,2).Value = "None"
$ws.Cells.Item(54,3).Value = 0
$ws.Cells.Item(54,4).Value = "None"
$ws.Cells.Item(55,1).Value = "Penn Foster"
$ws.Cells.Item(55,2).Value = "None"
$ws.Cells.Item(55,3).Value = 0
$ws.Cells.Item(55,4).Value = "None"
$ws.Cells.Item(56,1).Value = "Pfizer Inc."
$ws.Cells.Item(56,2).Value = "None"
$ws.Cells.Item(56,3).Value = 0
$ws.Cells.Item(56,4).Value = "None"
$ws.Cells.Item(57,1).Value = "Queen Rania Foundation for Education and Development"
$ws.Cells.Item(57,2).Value = "None"
$ws.Cells.Item(57,3).Value = 0
$ws.Cells.Item(57,4).Value = "None"
$ws.Cells.Item(58,1).Value = "RISE"
$ws.Cells.Item(58,2).Value = "None"
$ws.Cells.Item(58,3).Value = 0
$ws.Cells.Item(58,4).Value = "None"
$ws.Cells.Item(59,1).Value = "Save the Children"
$ws.Cells.Item(59,2).Value = "None"
$ws.Cells.Item(59,3).Value = 0
$ws.Cells.Item(59,4).Value = "None"
$ws.Cells.Item(60,1).Value = "Seed Global Health"
$ws.Cells.Item(60,2).Value = "None"
$ws.Cells.Item(60,3).Value = 0
$ws.Cells.Item(60,4).Value = "None"
$ws.Cells.Item(61,1).Value = "Someone Else's Child Foundation"
$ws.Cells.Item(61,2).Value = "None"
$ws.Cells.Item(61,3).Value = 0
$ws.Cells.Item(61,4).Value = "None"
$ws.Cells.Item(62,1).Value = "Soronko Solutions"
$ws.Cells.Item(62,2).Value = "None"
$ws.Cells.Item(62,3).Value = 0
$ws.Cells.Item(62,4).Value = "None"
$ws.Cells.Item(63,1).Value = "Sresta Natural Bioproducts Pvt Ltd"
$ws.Cells.Item(63,2).Value = "None"
$ws.Cells.Item(63,3).Value = 0
$ws.Cells.Item(63,4).Value = "None"
$ws.Cells.Item(64,1).Value = "Stand Together"
$ws.Cells.Item(64,2).Value = "None"
$ws.Cells.Item(64,3).Value = 0
$ws.Cells.Item(64,4).Value = "None"
$ws.Cells.Item(65,1).Value = "Strada Education Network"
$ws.Cells.Item(65,2).Value = "None"
$ws.Cells.Item(65,3).Value = 0
$ws.Cells.Item(65,4).Value = "None"
$ws.Cells.Item(66,1).Value = "Tecnológico de Monterrey"
$ws.Cells.Item(66,2).Value = "None"
$ws.Cells.Item(66,3).Value = 0
$ws.Cells.Item(66,4).Value = "None"
$ws.Cells.Item(67,1).Value = "TGR Foundation"
$ws.Cells.Item(67,2).Value = "None"
$ws.Cells.Item(67,3).Value = 0
$ws.Cells.Item(67,4).Value = "None"
$ws.Cells.Item(68,1).Value = "The Hague Business Agency"
$ws.Cells.Item(68,2).Value = "None"
$ws.Cells.Item(68,3).Value = 0
$ws.Cells.Item(68,4).Value = "None"
$ws.Cells.Item(69,1).Value = "The Kamath Family Foundation"
$ws.Cells.Item(69,2).Value = "None,AHSA Platform"
$ws.Cells.Item(69,3).Value = 1
$ws.Cells.Item(69,4).Value = "None"
$ws.Cells.Item(70,1).Value = "The Nature Conservancy"
$ws.Cells.Item(70,2).Value = "None"
$ws.Cells.Item(70,3).Value = 0
$ws.Cells.Item(70,4).Value = "None"
$ws.Cells.Item(71,1).Value = "The Pershing Square Foundation"
$ws.Cells.Item(71,2).Value = "None"
$ws.Cells.Item(71,3).Value = 0
$ws.Cells.Item(71,4).Value = "None"
$ws.Cells.Item(72,1).Value = "Twilio.org"
$ws.Cells.Item(72,2).Value = "None"
$ws.Cells.Item(72,3).Value = 0
$ws.Cells.Item(72,4).Value = "None"
$ws.Cells.Item(73,1).Value = "Uber"
$ws.Cells.Item(73,2).Value = "None"
$ws.Cells.Item(73,3).Value = 0
$ws.Cells.Item(73,4).Value = "None"
$ws.Cells.Item(74,1).Value = "Ultranauts Inc"
$ws.Cells.Item(74,2).Value = "None"
$ws.Cells.Item(74,3).Value = 0
$ws.Cells.Item(74,4).Value = "None"
$ws.Cells.Item(75,1).Value = "UN Women"
$ws.Cells.Item(75,2).Value = "None"
$ws.Cells.Item(75,3).Value = 0
$ws.Cells.Item(75,4).Value = "None"
$ws.Cells.Item(76,1).Value = "Usizo Advisory Solutions"
$ws.Cells.Item(76,2).Value = "None,AHSA Platform"
$ws.Cells.Item(76,3).Value = 1
$ws.Cells.Item(76,4).Value = "None"
$ws.Cells.Item(77,1).Value = "Women’s WorldWide Web (W4)"
$ws.Cells.Item(77,2).Value = "None"
$ws.Cells.Item(77,3).Value = 0
$ws.Cells.Item(77,4).Value = "None"
$ws.Cells.Item(78,1).Value = "Xprize"
$ws.Cells.Item(78,2).Value = "None"
$ws.Cells.Item(78,3).Value = 0
$ws.Cells.Item(78,4).Value = "None"
$ws.Cells.Item(79,1).Value = "YUM Brands"
$ws.Cells.Item(79,2).Value = "None"
$ws.Cells.Item(79,3).Value = 0
$ws.Cells.Item(79,4).Value = "None"

# Restore original active sheet/selection
$wb.Worksheets.Item(1).Activate()
